$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
  $r = $ws.Range($addr)
  $r.NumberFormat = "@"
  $r.Value = $val
  $r.Style = "Normal"
}

Set-TextCell 'A2' 'India'
Set-TextCell 'B2' 'Hyderabad'
Set-TextCell 'C2' 'Hyderabad City'
Set-TextCell 'D2' 'Saroor Nagar'
Set-TextCell 'E2' 'Aroma Family Restaurant'
Set-TextCell 'F2' 'Aroma Family Restaurant'
Set-TextCell 'G2' '4.1'
Set-TextCell 'H2' '209'
Set-TextCell 'I2' '3.6'
Set-TextCell 'J2' '14.5K'
Set-TextCell 'K2' '[''Hyderabadi'', ''North Indian'', ''Chinese'', ''Biryani'', ''Shawarma'']'
Set-TextCell 'L2' 'Saroor Nagar, Hyderabad'
Set-TextCell 'M2' '[''Hyderabadi'', ''North Indian'', ''Chinese'', ''Biryani'', ''Shawarma'']'
Set-TextCell 'N2' '[]'
Set-TextCell 'O2' '[''₹700 for two people (approx.)'']'
Set-TextCell 'P2' '[''Home Delivery'', ''Takeaway Available'', ''Indoor Seating'', ''Family Friendly'']'

Set-TextCell 'A3' 'India'
Set-TextCell 'B3' 'Hyderabad'
Set-TextCell 'C3' 'Hyderabad City'
Set-TextCell 'D3' 'Nampally'
Set-TextCell 'E3' 'Nimrah Restaurant'
Set-TextCell 'F3' 'Nimrah Restaurant'
Set-TextCell 'G3' '3.9'
Set-TextCell 'H3' '71'
Set-TextCell 'I3' '4.1'
Set-TextCell 'J3' '461'
Set-TextCell 'K3' '[''Bakery'', ''Fast Food'']'
Set-TextCell 'L3' 'Nampally, Hyderabad'
Set-TextCell 'M3' '[''Bakery'', ''Fast Food'']'
Set-TextCell 'N3' '[]'
Set-TextCell 'O3' '[''₹200 for two people (approx.)'']'
Set-TextCell 'P3' '[''Breakfast'', ''Home Delivery'', ''Takeaway Available'', ''Desserts and Bakes'']'

Set-TextCell 'A4' 'India'
Set-TextCell 'B4' 'Hyderabad'
Set-TextCell 'C4' 'Hyderabad City'
Set-TextCell 'D4' 'RTC X roads'
Set-TextCell 'E4' 'Sahara Bakers'
Set-TextCell 'F4' 'Sahara Bakers'
Set-TextCell 'G4' '3.8'
Set-TextCell 'H4' '723'
Set-TextCell 'I4' '3.5'
Set-TextCell 'J4' '69.1K'
Set-TextCell 'K4' '[''Chinese'', ''Bakery'', ''Sichuan'', ''Pizza'', ''Burger'', ''Fast Food'', ''Desserts'']'
Set-TextCell 'L4' 'RTC X roads, Hyderabad'
Set-TextCell 'M4' '[''Chinese'', ''Bakery'', ''Sichuan'', ''Pizza'', ''Burger'', ''Fast Food'', ''Desserts'']'
Set-TextCell 'N4' '[''Coffee and Doughnuts, Yummy Cake, Choco Chip Cake, Fruit Biscuits, Strawberry Cake, Butterscotch Pastry'']'
Set-TextCell 'O4' '[''₹300 for two people (approx.)'']'
Set-TextCell 'P4' '[''Home Delivery'', ''Takeaway Available'', ''Desserts and Bakes'', ''Indoor Seating'']'

Set-TextCell 'A5' 'India'
Set-TextCell 'B5' 'Hyderabad'
Set-TextCell 'C5' 'Hyderabad City'
Set-TextCell 'D5' 'Himayath Nagar'
Set-TextCell 'E5' 'McDonald''s'
Set-TextCell 'F5' 'McDonald''s'
Set-TextCell 'G5' '4.0'
Set-TextCell 'H5' '881'
Set-TextCell 'I5' '4.1'
Set-TextCell 'J5' '17.9K'
Set-TextCell 'K5' '[''Burger'', ''Fast Food'']'
Set-TextCell 'L5' 'Himayath Nagar, Hyderabad'
Set-TextCell 'M5' '[''Burger'', ''Fast Food'']'
Set-TextCell 'N5' '[''Mc Spicy Chicken, Chocolate Icecreams, Burgers, French Fries'']'
Set-TextCell 'O5' '[''₹500 for two people (approx.)'']'
Set-TextCell 'P5' '[''Home Delivery'', ''Takeaway Available'', ''Indoor Seating'']'

Set-TextCell 'A6' 'India'
Set-TextCell 'B6' 'Hyderabad'
Set-TextCell 'C6' 'Hyderabad City'
Set-TextCell 'D6' 'Lakdikapul'
Set-TextCell 'E6' 'Chicha''s'
Set-TextCell 'F6' 'Chicha''s'
Set-TextCell 'G6' '4.0'
Set-TextCell 'H6' '1,902'
Set-TextCell 'I6' '3.9'
Set-TextCell 'J6' '306'
Set-TextCell 'K6' '[''Kebab'', ''Rolls'', ''Chinese'', ''Biryani'', ''Desserts'', ''Beverages'']'
Set-TextCell 'L6' 'Lakdikapul, Hyderabad'
Set-TextCell 'M6' '[''Kebab'', ''Rolls'', ''Chinese'', ''Biryani'', ''Desserts'', ''Beverages'']'
Set-TextCell 'N6' '[''Lamb Briyani, Keema Masala, Bheja Fry, Red Chicken, Pathar Ka Ghost, Authentic Hyderabadi Food'']'
Set-TextCell 'O6' '[''₹1,000 for two people (approx.)'']'
Set-TextCell 'P6' '[''Home Delivery'', ''Takeaway Available'', ''Outdoor Seating'', ''Family Friendly'', ''Indoor Seating'', ''Desserts and Bakes'']'

Set-TextCell 'A7' 'India'
Set-TextCell 'B7' 'Hyderabad'
Set-TextCell 'C7' 'Hyderabad City'
Set-TextCell 'D7' 'Charminar'
Set-TextCell 'E7' 'Shah Ghouse Hotel & Restaurant'
Set-TextCell 'F7' 'Shah Ghouse Hotel & Restaurant'
Set-TextCell 'G7' '4.2'
Set-TextCell 'H7' '3,521'
Set-TextCell 'I7' '4.2'
Set-TextCell 'J7' '120.4K'
Set-TextCell 'K7' '[''North Indian'', ''Mughlai'', ''Chinese'', ''Mandi'', ''Biryani'', ''Shawarma'', ''Desserts'']'
Set-TextCell 'L7' 'Charminar, Hyderabad'
Set-TextCell 'M7' '[''North Indian'', ''Mughlai'', ''Chinese'', ''Mandi'', ''Biryani'', ''Shawarma'', ''Desserts'']'
Set-TextCell 'N7' '[''Hariyali Chicken, Authentic Hyderabadi Biryani, Mutton Haleem, Plain Rice, Chai, Tandoori Chicken'']'
Set-TextCell 'O7' '[''₹1,000 for two people (approx.)'']'
Set-TextCell 'P7' '[''Home Delivery'', ''Takeaway Available'', ''Indoor Seating'', ''Family Friendly'', ''Desserts and Bakes'']'

Set-TextCell 'A8' 'India'
Set-TextCell 'B8' 'Hyderabad'
Set-TextCell 'C8' 'Hyderabad City'
Set-TextCell 'D8' 'Lakdikapul'
Set-TextCell 'E8' 'Peshawar'
Set-TextCell 'F8' 'Peshawar'
Set-TextCell 'G8' '3.5'
Set-TextCell 'H8' '807'
Set-TextCell 'I8' '4.0'
Set-TextCell 'J8' '254'
Set-TextCell 'K8' '[''North Indian'', ''Seafood'', ''Kebab'', ''Chinese'']'
Set-TextCell 'L8' 'Lakdikapul, Hyderabad'
Set-TextCell 'M8' '[''North Indian'', ''Seafood'', ''Kebab'', ''Chinese'']'
Set-TextCell 'N8' '[]'
Set-TextCell 'O8' '[''₹1,200 for two people (approx.)'']'
Set-TextCell 'P8' '[''Home Delivery'', ''Takeaway Available'', ''Indoor Seating'', ''Family Friendly'']'

Set-TextCell 'A9' 'India'
Set-TextCell 'B9' 'Hyderabad'
Set-TextCell 'C9' 'Hyderabad City'
Set-TextCell 'D9' 'Abids'
Set-TextCell 'E9' 'Krupa Mess & Tiffins'
Set-TextCell 'F9' 'Krupa Mess & Tiffins'
Set-TextCell 'G9' '3.4'
Set-TextCell 'H9' '615'
Set-TextCell 'I9' '3.6'
Set-TextCell 'J9' '40.3K'
Set-TextCell 'K9' '[''South Indian'', ''Chinese'', ''North Indian'', ''Sichuan'', ''Pizza'']'
Set-TextCell 'L9' 'Abids, Hyderabad'
Set-TextCell 'M9' '[''South Indian'', ''Chinese'', ''North Indian'', ''Sichuan'', ''Pizza'']'
Set-TextCell 'N9' '[]'
Set-TextCell 'O9' '[''₹300 for two people (approx.)'']'
Set-TextCell 'P9' '[''Breakfast'', ''Home Delivery'', ''Takeaway Available'', ''Indoor Seating'']'

Set-TextCell 'A10' 'India'
Set-TextCell 'B10' 'Hyderabad'
Set-TextCell 'C10' 'Hyderabad City'
Set-TextCell 'D10' 'Narayanguda'
Set-TextCell 'E10' 'Mehfil'
Set-TextCell 'F10' 'Mehfil'
Set-TextCell 'G10' '4.1'
Set-TextCell 'H10' '4,685'
Set-TextCell 'I10' '4.0'
Set-TextCell 'J10' '185.3K'
Set-TextCell 'K10' '[''Kebab'', ''Biryani'']'
Set-TextCell 'L10' 'Narayanguda, Hyderabad'
Set-TextCell 'M10' '[''Kebab'', ''Biryani'']'
Set-TextCell 'N10' '[''Jumbo Chicken Biryani, Panneer Butter Masala, Naan, Tea'']'
Set-TextCell 'O10' '[''₹700 for two people (approx.)'']'
Set-TextCell 'P10' '[''Home Delivery'', ''Takeaway Available'', ''Family Friendly'', ''Indoor Seating'']'
